$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old content first (A1:A3 had FUCK/YEAH/NO)
$ws.Range("A1:A3").ClearContents()

# New header row (bold): PRODUCT, PRICE, THEME, BUTTON TYPE
$ws.Range("A1").Value = "PRODUCT"
$ws.Range("B1").Value = "PRICE"
$ws.Range("C1").Value = "THEME"
$ws.Range("D1").Value = "BUTTON TYPE"
$ws.Range("A1:D1").Font.Bold = $true

# New data row
$ws.Range("A2").Value = "test"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "20"
$ws.Range("B2").ClearFormats()

# Column widths matching the diff (closest achievable given Excel's pixel quantization)
$ws.Columns.Item(1).ColumnWidth = 8.833333333333334
$ws.Columns.Item(2).ColumnWidth = 5.166666666666667
$ws.Columns.Item(3).ColumnWidth = 6.666666666666667
$ws.Columns.Item(4).ColumnWidth = 12.5
